$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $escaped = $val -replace '"', '""'
    $rng.Formula = '="' + $escaped + '"'
    $rng.Copy() | Out-Null
    $rng.PasteSpecial([Microsoft.Office.Interop.Excel.XlPasteType]::xlPasteValues) | Out-Null
}

$ws.Range("D2").Value = '30.583.49'
$ws.Range("E2").Value = '  -0.41%  '

$ws.Range("D3").Value = '1.882.52'
$ws.Range("E3").Value = '  -0.30%  '

Set-TextValue "D4" '0.9997'
$ws.Range("E4").Value = '  -0.12%  '

Set-TextValue "D5" '246.38'
$ws.Range("E5").Value = '  -0.62%  '

Set-TextValue "D6" '0.9998'
$ws.Range("E6").Value = '  -0.07%  '

Set-TextValue "D7" '0.4728'
$ws.Range("E7").Value = '  -0.09%  '

Set-TextValue "D8" '0.2885'

Set-TextValue "D9" '0.06536'
$ws.Range("E9").Value = '  +0.12%  '

Set-TextValue "D10" '22.14'
$ws.Range("E10").Value = '  +0.54%  '

$ws.Range("B11").Value = 'Polygon'
$ws.Range("C11").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
Set-TextValue "D11" '0.7730'
$ws.Range("E11").Value = '  +5.17%  '

$ws.Range("B12").Value = 'Litecoin'
$ws.Range("C12").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
Set-TextValue "D12" '100.90'
$ws.Range("E12").Value = '  +4.40%  '

Set-TextValue "D13" '0.07827'
$ws.Range("E13").Value = '  +0.36%  '

$ws.Range("D14").Value = '1.881.79'
$ws.Range("E14").Value = '  -0.41%  '

Set-TextValue "D15" '5.251'
$ws.Range("E15").Value = '  +0.02%  '

Set-TextValue "D16" '285.46'
$ws.Range("E16").Value = '  +0.69%  '

$ws.Range("D17").Value = '30.544.29'
$ws.Range("E17").Value = '  -0.51%  '

Set-TextValue "D18" '13.20'
$ws.Range("E18").Value = '  -0.25%  '

Set-TextValue "D19" '0.000007525'
$ws.Range("E19").Value = '  -0.03%  '

$ws.Range("E20").Value = '  -0.04%  '

$ws.Range("D21").Value = '2.125.76'
$ws.Range("E21").Value = '  -0.68%  '

Set-TextValue "D22" '5.370'
$ws.Range("E22").Value = '  +1.06%  '

Set-TextValue "D23" '0.9994'
$ws.Range("E23").Value = '  -0.10%  '

Set-TextValue "D24" '6.405'
$ws.Range("E24").Value = '  +2.55%  '

Set-TextValue "D25" '9.137'
$ws.Range("E25").Value = '  -0.90%  '

Set-TextValue "D26" '162.93'
$ws.Range("E26").Value = '  -0.95%  '

Set-TextValue "D27" '19.09'
$ws.Range("E27").Value = '  +0.77%  '

Set-TextValue "D28" '1.916'
$ws.Range("E28").Value = '  -0.16%  '

Set-TextValue "D29" '0.09703'
$ws.Range("E29").Value = '  -0.29%  '

Set-TextValue "D30" '1.329'
$ws.Range("E30").Value = '  -0.75%  '

$ws.Range("E31").Value = '  +0.94%  '

Set-TextValue "D32" '4.267'
$ws.Range("E32").Value = '  -0.71%  '

Set-TextValue "D33" '4.199'
$ws.Range("E33").Value = '  -0.06%  '

Set-TextValue "D34" '0.04848'
$ws.Range("E34").Value = '  -0.37%  '

Set-TextValue "D35" '1.130'
$ws.Range("E35").Value = '  +0.28%  '

Set-TextValue "D36" '0.6967'
$ws.Range("E36").Value = '  -0.18%  '

Set-TextValue "D37" '2.752'
$ws.Range("E37").Value = '  +0.98%  '

Set-TextValue "D38" '0.01916'
$ws.Range("E38").Value = '  +1.25%  '

$ws.Range("E39").Value = '  +2.70%  '

Set-TextValue "D40" '76.60'
$ws.Range("E40").Value = '  +0.65%  '

Set-TextValue "D41" '6.289'
$ws.Range("E41").Value = '  -1.28%  '

Set-TextValue "D42" '1.980'
$ws.Range("E42").Value = '  -1.08%  '

Set-TextValue "D43" '0.4266'
$ws.Range("E43").Value = '  +0.36%  '

Set-TextValue "D44" '0.9995'
$ws.Range("E44").Value = '  -0.14%  '

Set-TextValue "D45" '0.8297'
$ws.Range("E45").Value = '  -0.86%  '

Set-TextValue "D46" '101.52'
$ws.Range("E46").Value = '  +0.02%  '

Set-TextValue "D47" '9.795'
$ws.Range("E47").Value = '  +2.81%  '

Set-TextValue "D48" '7.039'
$ws.Range("E48").Value = '  +0.21%  '

Set-TextValue "D49" '35.14'
$ws.Range("E49").Value = '  -1.53%  '

Set-TextValue "D50" '896.92'
$ws.Range("E50").Value = '  -2.11%  '

Set-TextValue "D51" '0.05762'
$ws.Range("E51").Value = '  +0.12%  '
